$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "27.683.55"
$ws.Cells.Item(2, 5).Value = "  +5.99%  "
$ws.Cells.Item(3, 4).Value = "1.734.93"
$ws.Cells.Item(3, 5).Value = "  +4.86%  "
$ws.Cells.Item(4, 5).Value = "  +0.01%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "227.70"
$ws.Cells.Item(5, 4).NumberFormat = "General"
$ws.Cells.Item(5, 5).Value = "  +4.05%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "0.5452"
$ws.Cells.Item(6, 4).NumberFormat = "General"
$ws.Cells.Item(6, 5).Value = "  +3.70%  "
$ws.Cells.Item(7, 5).Value = "  -0.02%  "
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.2745"
$ws.Cells.Item(8, 4).NumberFormat = "General"
$ws.Cells.Item(8, 5).Value = "  +2.24%  "
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.06716"
$ws.Cells.Item(9, 4).NumberFormat = "General"
$ws.Cells.Item(9, 5).Value = "  +5.40%  "
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "21.94"
$ws.Cells.Item(10, 4).NumberFormat = "General"
$ws.Cells.Item(10, 5).Value = "  +6.66%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.07790"
$ws.Cells.Item(11, 4).NumberFormat = "General"
$ws.Cells.Item(11, 5).Value = "  +1.23%  "
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "4.696"
$ws.Cells.Item(12, 4).NumberFormat = "General"
$ws.Cells.Item(12, 5).Value = "  +1.79%  "
$ws.Cells.Item(13, 4).Value = "1.755.95"
$ws.Cells.Item(13, 5).Value = "  +1.92%  "
$ws.Cells.Item(14, 4).Value = "1.974.44"
$ws.Cells.Item(14, 5).Value = "  +4.81%  "
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "0.5989"
$ws.Cells.Item(15, 4).NumberFormat = "General"
$ws.Cells.Item(15, 5).Value = "  +6.37%  "
$ws.Cells.Item(16, 4).Value = "0.0₅8423"
$ws.Cells.Item(16, 5).Value = "  +1.94%  "
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "69.23"
$ws.Cells.Item(17, 4).NumberFormat = "General"
$ws.Cells.Item(17, 5).Value = "  +5.45%  "
$ws.Cells.Item(18, 4).Value = "27.696.88"
$ws.Cells.Item(18, 5).Value = "  +6.08%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "228.89"
$ws.Cells.Item(19, 4).NumberFormat = "General"
$ws.Cells.Item(19, 5).Value = "  +20.36%  "
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "4.835"
$ws.Cells.Item(20, 4).NumberFormat = "General"
$ws.Cells.Item(20, 5).Value = "  +2.98%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "1.003"
$ws.Cells.Item(21, 4).NumberFormat = "General"
$ws.Cells.Item(21, 5).Value = "  -0.04%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "10.90"
$ws.Cells.Item(22, 4).NumberFormat = "General"
$ws.Cells.Item(22, 5).Value = "  +5.29%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "6.226"
$ws.Cells.Item(23, 4).NumberFormat = "General"
$ws.Cells.Item(23, 5).Value = "  +3.79%  "
$ws.Cells.Item(24, 5).Value = "  +0.01%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "148.40"
$ws.Cells.Item(25, 4).NumberFormat = "General"
$ws.Cells.Item(25, 5).Value = "  +1.04%  "
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "1.731"
$ws.Cells.Item(26, 4).NumberFormat = "General"
$ws.Cells.Item(26, 5).Value = "  +13.71%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "0.1248"
$ws.Cells.Item(27, 4).NumberFormat = "General"
$ws.Cells.Item(27, 5).Value = "  +3.86%  "
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "7.471"
$ws.Cells.Item(28, 4).NumberFormat = "General"
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "17.14"
$ws.Cells.Item(29, 4).NumberFormat = "General"
$ws.Cells.Item(29, 5).Value = "  +7.18%  "
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "0.05711"
$ws.Cells.Item(30, 4).NumberFormat = "General"
$ws.Cells.Item(30, 5).Value = "  +1.12%  "
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "1.313"
$ws.Cells.Item(31, 4).NumberFormat = "General"
$ws.Cells.Item(31, 5).Value = "  +2.85%  "
$ws.Cells.Item(32, 5).Value = "  +5.82%  "
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "3.520"
$ws.Cells.Item(33, 4).NumberFormat = "General"
$ws.Cells.Item(33, 5).Value = "  +4.02%  "
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "1.690"
$ws.Cells.Item(34, 4).NumberFormat = "General"
$ws.Cells.Item(34, 5).Value = "  +6.84%  "
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "0.9763"
$ws.Cells.Item(35, 4).NumberFormat = "General"
$ws.Cells.Item(35, 5).Value = "  +2.91%  "
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "2.860"
$ws.Cells.Item(36, 4).NumberFormat = "General"
$ws.Cells.Item(36, 5).Value = "  +2.19%  "
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "2.440"
$ws.Cells.Item(37, 4).NumberFormat = "General"
$ws.Cells.Item(37, 5).Value = "  +1.24%  "
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.5986"
$ws.Cells.Item(38, 4).NumberFormat = "General"
$ws.Cells.Item(38, 5).Value = "  +3.55%  "
$ws.Cells.Item(39, 5).Value = "  +4.71%  "
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "5.958"
$ws.Cells.Item(40, 4).NumberFormat = "General"
$ws.Cells.Item(40, 5).Value = "  -0.26%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.8506"
$ws.Cells.Item(41, 4).NumberFormat = "General"
$ws.Cells.Item(41, 5).Value = "  +1.72%  "
$ws.Cells.Item(42, 4).Value = "1.050.17"
$ws.Cells.Item(42, 5).Value = "  +2.70%  "
$ws.Cells.Item(43, 5).Value = "  +0.05%  "
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "101.59"
$ws.Cells.Item(44, 4).NumberFormat = "General"
$ws.Cells.Item(44, 5).Value = "  +0.20%  "
$ws.Cells.Item(45, 4).Value = "1.879.50"
$ws.Cells.Item(45, 5).Value = "  +4.78%  "
$ws.Cells.Item(46, 4).Value = "0.0₈116"
$ws.Cells.Item(46, 5).Value = "  +9.81%  "
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "59.63"
$ws.Cells.Item(47, 4).NumberFormat = "General"
$ws.Cells.Item(47, 5).Value = "  +2.12%  "
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "8.284"
$ws.Cells.Item(48, 4).NumberFormat = "General"
$ws.Cells.Item(48, 5).Value = "  +3.20%  "
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "0.4437"
$ws.Cells.Item(49, 4).NumberFormat = "General"
$ws.Cells.Item(49, 5).Value = "  +2.18%  "
$ws.Cells.Item(50, 2).Value = "Cronos"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "0.05337"
$ws.Cells.Item(50, 4).NumberFormat = "General"
$ws.Cells.Item(50, 5).Value = "  -0.15%  "
$ws.Cells.Item(51, 2).Value = "Frax"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "1.002"
$ws.Cells.Item(51, 4).NumberFormat = "General"
$ws.Cells.Item(51, 5).Value = "  -0.15%  "
